$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '30.008.21'
$ws.Range('D2').ClearFormats()
$ws.Range('E2').Value = '  -0.18%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.878.46'
$ws.Range('D3').ClearFormats()
$ws.Range('E3').Value = '  -0.40%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9999'
$ws.Range('D4').ClearFormats()
$ws.Range('E4').Value = '  +0.24%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '241.40'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -3.16%  '
$ws.Range('E6').Value = '  +0.19%  '
$ws.Range('E7').Value = '  -0.70%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '44.53'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  -2.62%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.2903'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  +1.35%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.06585'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  +0.33%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '1.879.02'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  +0.02%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '16.67'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  -3.14%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.07184'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  -0.69%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.6615'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  -1.26%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '85.95'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  +1.06%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '4.840'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  +0.31%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '29.973.04'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  -0.28%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.000007888'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  +4.68%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.9996'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  +0.08%  '
$ws.Range('E20').Value = '  -1.75%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '2.119.46'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  -0.10%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '1.000'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  +0.33%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '4.747'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  -0.60%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '5.582'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  +0.80%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '9.075'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  +0.53%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '149.99'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  +3.15%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '134.19'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  -0.38%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '16.70'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  -0.33%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.899'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  -2.75%  '
$ws.Range('E30').Value = '  +0.23%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.154'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  -0.94%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.08695'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  +0.42%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.930'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  +0.34%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.05012'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  -0.67%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.7027'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  +1.13%  '
$ws.Range('E36').Value = '  -4.08%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.654'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  -1.32%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.687'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  -2.33%  '
$ws.Range('E39').Value = '  -5.44%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.01696'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  +3.25%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.9310'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  -3.37%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '5.956'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  -2.03%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.9989'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  +0.03%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.4177'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  -1.07%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '101.44'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  -3.02%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '7.399'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  -0.64%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.1255'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  -0.25%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.05663'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  +0.33%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '32.37'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  -0.27%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '55.70'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  +1.34%  '
$ws.Range('B51').Value = 'EnergySwap'
$ws.Range('C51').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '8.096'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  -1.99%  '
